$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: lMS9P2yB | Bala vs Connahs Q.
$rowData = @("lMS9P2yB", "15/10/2024", "15:45", "WALES - CYMRU PREMIER", "Bala", "Connahs Q.", 2.37, 3.1, 2.82, 2.95, 2.1, 3.35, 1.07, 7, 1.32, 3.1, 1.95, 1.78, 1.38, 2.8, 1.72, 2.02, 7.9, 11.75, 9.25, 25, 20, 29, 7, 6.1, 13, 60, 9.25, 15, 10.25, 35, 24, 32, 450, 4.4, 12.5, 19, 50, 80, 250, 2.8, 6.6, 55, 4.9, 15.5, 21, 70, 90, 250, 51, 51)
for ($i = 0; $i -lt $rowData.Length; $i++) {
    $ws.Cells.Item(4, $i + 1).Value = $rowData[$i]
}

# Row 5: rLLINt6N | Briton Ferry vs Haverfordwest
$rowData = @("rLLINt6N", "15/10/2024", "15:45", "WALES - CYMRU PREMIER", "Briton Ferry", "Haverfordwest", 6.9, 4.25, 1.39, 6.6, 2.3, 1.91, 1.05, 7.8, 1.26, 3.45, 1.8, 1.95, 1.37, 2.85, 2.05, 1.7, 17, 45, 23, 175, 90, 80, 7.8, 8.5, 21, 110, 6.3, 6.2, 8.5, 8.75, 12, 30, 201, 8.25, 45, 45, 300, 350, 51, 2.85, 8.75, 90, 3.15, 6.4, 18, 18.5, 55, 300, 51, 51)
for ($i = 0; $i -lt $rowData.Length; $i++) {
    $ws.Cells.Item(5, $i + 1).Value = $rowData[$i]
}

# Row 6: vi7lIKTp | Caernarfon vs Aberystwyth
$rowData = @("vi7lIKTp", "15/10/2024", "15:45", "WALES - CYMRU PREMIER", "Caernarfon", "Aberystwyth", 1.31, 4.85, 7.6, 1.72, 2.6, 6.8, 1.03, 9.25, 1.16, 4.65, 1.5, 2.42, 1.27, 3.45, 1.78, 1.93, 8.75, 7.4, 8.5, 8.75, 10.25, 23, 9.25, 10.25, 19, 75, 24, 55, 24, 200, 80, 65, 500, 3.35, 5.8, 14, 14.5, 35, 150, 3.45, 8, 65, 9, 45, 40, 300, 250, 450, 51, 51)
for ($i = 0; $i -lt $rowData.Length; $i++) {
    $ws.Cells.Item(6, $i + 1).Value = $rowData[$i]
}

# Row 7: fakytOLi | Cardiff Metropolitan vs Penybont
$rowData = @("fakytOLi", "15/10/2024", "15:45", "WALES - CYMRU PREMIER", "Cardiff Metropolitan", "Penybont", 4, 3.25, 1.85, 4.5, 2.12, 2.42, 1.08, 6.6, 1.36, 2.9, 2.07, 1.7, 1.42, 2.67, 1.9, 1.8, 10, 21, 14, 65, 45, 50, 6.6, 6.4, 16, 90, 6.4, 8.25, 8.25, 15.5, 15.5, 30, 800, 5.9, 24, 30, 150, 175, 400, 2.67, 7.4, 70, 3.7, 9.25, 18.5, 32, 65, 250, 51, 51)
for ($i = 0; $i -lt $rowData.Length; $i++) {
    $ws.Cells.Item(7, $i + 1).Value = $rowData[$i]
}

# Row 8: 2Phqv2j4 | Newtown vs Barry
$rowData = @("2Phqv2j4", "15/10/2024", "15:45", "WALES - CYMRU PREMIER", "Newtown", "Barry", 2.3, 3.25, 2.82, 2.87, 2.12, 3.4, 1.05, 7.7, 1.25, 3.55, 1.78, 1.98, 1.37, 2.85, 1.6, 2.2, 9.25, 12.5, 9, 25, 17.5, 24, 7.7, 6.4, 12, 45, 10.25, 16, 10.25, 35, 23, 27, 300, 4.4, 12, 18.5, 50, 75, 200, 2.85, 6.6, 50, 5, 15.5, 21, 70, 90, 250, 51, 51)
for ($i = 0; $i -lt $rowData.Length; $i++) {
    $ws.Cells.Item(8, $i + 1).Value = $rowData[$i]
}
